$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Red_TestData")

# --- G3: severity changed from "Critical" to "High " ---
$ws.Range("G3").Value = "High "

# --- Copy row 3's formatting down into the three new rows (4,5,6) so the ---
# --- new test-data rows inherit the same look (borders/fonts/number fmt) ---
$ws.Range("A3:R3").Copy()
$ws.Range("A4:R6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 4 : John Rease / Testing / proj-003 ---
$ws.Range("A4").Value = "June 13,2023"
$ws.Range("B4").Value = "John Rease"
$ws.Range("C4").Value = "Testing"
$ws.Range("D4").Value = "https://jira.com/browse/proj-003"
$ws.Range("E4").Value = "M109"
$ws.Range("F4").Value = "Yes "
$ws.Range("G4").Value = "Critical"
$ws.Range("H4").Value = "High "
$ws.Range("I4").Value = "High "
$ws.Range("J4").Value = "N/A"
$ws.Range("K4").Value = "N/A"
$ws.Range("L4").Value = "Critical"
$ws.Range("M4").Value = "N/A"
$ws.Range("N4").Value = "N/A"
$ws.Range("O4").Value = "Yes "
$ws.Range("P4").Value = 44706
$ws.Range("Q4").Value = "No"
$ws.Range("R4").Value = "No"
$ws.Hyperlinks.Add($ws.Range("D4"), "https://jira.com/browse/proj-003") | Out-Null

# --- Row 5 : Harold Finch / PO Review / proj-004 ---
$ws.Range("A5").Value = "June 13,2023"
$ws.Range("B5").Value = "Harold Finch"
$ws.Range("C5").Value = "PO Review"
$ws.Range("D5").Value = "https://jira.com/browse/proj-004"
$ws.Range("E5").Value = "M110"
$ws.Range("F5").Value = "No "
$ws.Range("G5").Value = "Low "
$ws.Range("H5").Value = "High "
$ws.Range("I5").Value = "High "
$ws.Range("J5").Value = "N/A"
$ws.Range("K5").Value = "N/A"
$ws.Range("L5").Value = "Critical"
$ws.Range("M5").Value = "N/A"
$ws.Range("N5").Value = "N/A"
$ws.Range("O5").Value = "Yes "
$ws.Range("P5").Value = 44706
$ws.Range("Q5").Value = "No"
$ws.Range("R5").Value = "No"
$ws.Hyperlinks.Add($ws.Range("D5"), "https://jira.com/browse/proj-004") | Out-Null

# --- Row 6 : Tony Wei / Sign Off / proj-005 ---
$ws.Range("A6").Value = "June 13,2023"
$ws.Range("B6").Value = "Tony Wei"
$ws.Range("C6").Value = "Sign Off"
$ws.Range("D6").Value = "https://jira.com/browse/proj-005"
$ws.Range("E6").Value = "M110"
$ws.Range("F6").Value = "Yes "
$ws.Range("G6").Value = "High "
$ws.Range("H6").Value = "N/A"
$ws.Range("I6").Value = "Medium"
$ws.Range("J6").Value = "N/A"
$ws.Range("K6").Value = "N/A"
$ws.Range("L6").Value = "Critical"
$ws.Range("M6").Value = "N/A"
$ws.Range("N6").Value = "N/A"
$ws.Range("O6").Value = "Yes "
$ws.Range("P6").Value = 44706
$ws.Range("Q6").Value = "No"
$ws.Range("R6").Value = "No"
$ws.Hyperlinks.Add($ws.Range("D6"), "https://jira.com/browse/proj-005") | Out-Null

# --- Update the active selection to match the author's final cursor position ---
$ws.Activate()
$ws.Range("I10").Select()
